$wb = $excel.ActiveWorkbook

function Set-MethodRow($ws, $row, $name, $c, $d) {
    $ws.Range("A$row").Value = $name
    $ws.Range("C$row").Value = $c
    if ($null -ne $d) {
        $ws.Range("D$row").Value = $d
    }
}

$ws = $wb.Worksheets.Item("FUSELAGE")
Set-MethodRow $ws 8 "ROSKAM" 14190.0 125.66337898688195
Set-MethodRow $ws 9 "NICOLAI_1984" 10196.0 62.14685075054604
Set-MethodRow $ws 10 "RAYMER" 6416.0 2.0335616335330875
Set-MethodRow $ws 11 "SADRAEY" 6396.0 1.7155019027552412
Set-MethodRow $ws 12 "JENKINSON" 21031.0 234.45570989944426
Set-MethodRow $ws 13 "KROO" 7092.0 12.783980533824291
Set-MethodRow $ws 14 "TORENBEEK_1976" 10802.0 71.78406059311477
Set-MethodRow $ws 15 "TORENBEEK_2013" 7891.0 25.49046677839925

$ws = $wb.Worksheets.Item("WING")
Set-MethodRow $ws 8 "RAYMER" 8327.0 24.928461235241798
Set-MethodRow $ws 9 "KROO" 7483.0 12.266083274086032
Set-MethodRow $ws 11 "TORENBEEK_2013" 6082.0 -8.752864028732963

$ws = $wb.Worksheets.Item("HORIZONTAL TAIL")
Set-MethodRow $ws 8 "NICOLAI_2013" 394.0 -45.51498524936025
Set-MethodRow $ws 9 "ROSKAM" 1523.0 110.61085651072167
Set-MethodRow $ws 10 "RAYMER" 502.0 -30.58000658674834
Set-MethodRow $ws 11 "SADRAEY" 1040.0 43.81831304737396
Set-MethodRow $ws 12 "HOWE" 1415.0 95.67587784810976
Set-MethodRow $ws 13 "JENKINSON" 700.0 -3.199212371959834
Set-MethodRow $ws 14 "TORENBEEK_1976" 52.0 -92.8090843476313
Set-MethodRow $ws 15 "KROO" 737.0 1.9174006883794317

$ws = $wb.Worksheets.Item("VERTICAL TAIL")
Set-MethodRow $ws 8 "ROSKAM" 1523.0 110.61085651072167
Set-MethodRow $ws 9 "RAYMER" 179.0 -75.24665573511544
Set-MethodRow $ws 10 "SADRAEY" 749.0 3.5768427620029772
Set-MethodRow $ws 11 "HOWE" 1145.0 58.338431191579986
Set-MethodRow $ws 12 "JENKINSON" 502.0 -30.58000658674834
Set-MethodRow $ws 13 "TORENBEEK_1976" 124.0 -82.85243190589003
Set-MethodRow $ws 14 "KROO" 485.0 -32.93088285771503

$ws = $wb.Worksheets.Item("NACELLES")
Set-MethodRow $ws 10 "ROSKAM" 687.0 15.003702654937063
Set-MethodRow $ws 12 "JENKINSON" 705.0 18.016900104411395
Set-MethodRow $ws 17 "ROSKAM" 687.0 15.003702654937063
Set-MethodRow $ws 19 "JENKINSON" 705.0 18.016900104411395

$ws = $wb.Worksheets.Item("POWER PLANT")
Set-MethodRow $ws 11 "KUNDU" 3265.0 25.116267589116617
Set-MethodRow $ws 12 "TORENBEEK_1976" 2954.0 13.19860779732021
Set-MethodRow $ws 13 "TORENBEEK_2013" 3458.0 32.51211434093882
Set-MethodRow $ws 18 "KUNDU" 3265.0 25.116267589116617
Set-MethodRow $ws 19 "TORENBEEK_1976" 2954.0 13.19860779732021
Set-MethodRow $ws 20 "TORENBEEK_2013" 3458.0 32.51211434093882

$ws = $wb.Worksheets.Item("LANDING GEARS")
Set-MethodRow $ws 9 "TORENBEEK_1976" 2499.2685173219097 -3.0589424681764026
Set-MethodRow $ws 11 "TORENBEEK_1976" 386.729549170154 $null
Set-MethodRow $ws 13 "TORENBEEK_1976" 2112.5389681517563 $null
